$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"

# Zero out billed amount / pricing figures (no-violation / reset scenario)
$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
